# Update bicicletas sheet: collapse TIPO/TOTAL columns into a single
# BICICLETAS column holding the totals, then drop the now-empty TOTAL column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the TOTAL values (column C) into column B, replacing the repeated
# "BICICLETAS" text values that used to live there.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $total
}

# New header for column B.
$ws.Cells.Item(1, 2).Value = "BICICLETAS"

# Remove the old column C (TOTAL) entirely, shifting nothing left of it.
$ws.Columns.Item(3).Delete()
